$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from the first (empty) paragraph.
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# 2. Footer: merge "Twitter: @" + "strategy_unit" into a single run/text,
#    removing the spell-check split.
$d.Content.Find.Execute("Twitter: @strategy_unit", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Twitter: @strategy_unit", 0)

# 3. Update the FootnoteText style (paragraph + run formatting).
$style = $d.Styles("FootnoteText")
$style.ParagraphFormat.SpaceAfter = 0
$style.ParagraphFormat.LineSpacingRule = 0
$style.Font.Name = "Segoe UI Light"
$style.Font.Size = 8
